$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 805
$ws.Range("I2").Value = 2162
$ws.Range("J2").Value = 8932
$ws.Range("L2").Value = 2448
$ws.Range("M2").Value = 153
$ws.Range("N2").Value = 1584
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 24
$ws.Range("Q2").Value = 12
$ws.Range("R2").Value = 109
$ws.Range("S2").Value = 928
$ws.Range("T2").Value = 1588
$ws.Range("U2").Value = 127
$ws.Range("V2").Value = 13895
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 13788
$ws.Range("Y2").Value = 30
$ws.Range("Z2").Value = 206
$ws.Range("AA2").Value = 83

$wb.Save()
